$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.049409018489615
$ws.Range("D2").Value = 1.061747717870855
$ws.Range("E2").Value = 1.057176091053682
$ws.Range("F2").Value = 1.069272477184303
$ws.Range("I2").Value = 1.0469462899829
$ws.Range("J2").Value = 1.054447620007442
$ws.Range("K2").Value = 1.064470937900029
$ws.Range("L2").Value = 1.059911778567524
$ws.Range("M2").Value = 1.071975427226106
$ws.Range("N2").Value = 1.021921283315917
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.05046414916293
$ws.Range("D3").Value = 1.062505536579119
$ws.Range("E3").Value = 1.058092608499733
$ws.Range("F3").Value = 1.07021964316242
$ws.Range("I3").Value = 1.047210371893458
$ws.Range("J3").Value = 1.055151501872561
$ws.Range("K3").Value = 1.065043581242572
$ws.Range("L3").Value = 1.060641819463084
$ws.Range("M3").Value = 1.072738406496183
$ws.Range("N3").Value = 1.022160700224695
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.051147288391048
$ws.Range("D4").Value = 1.062996140510159
$ws.Range("E4").Value = 1.058686337480131
$ws.Range("F4").Value = 1.070833229520407
$ws.Range("I4").Value = 1.047380233274471
$ws.Range("J4").Value = 1.055606764012561
$ws.Range("K4").Value = 1.065413698581937
$ws.Range("L4").Value = 1.061114252852189
$ws.Range("M4").Value = 1.073232182853129
$ws.Range("N4").Value = 1.022315423943418
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.051434575403526
$ws.Range("D5").Value = 1.063202447800588
$ws.Range("E5").Value = 1.058936103276986
$ws.Range("F5").Value = 1.071091349382238
$ws.Range("I5").Value = 1.047451399067704
$ws.Range("J5").Value = 1.055798108630453
$ws.Range("K5").Value = 1.065569194262764
$ws.Range("L5").Value = 1.061312874688046
$ws.Range("M5").Value = 1.073439784213006
$ws.Range("N5").Value = 1.02238042281606
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.051482817736634
$ws.Range("D6").Value = 1.063237091034894
$ws.Range("E6").Value = 1.058978049535441
$ws.Range("F6").Value = 1.071134698674227
$ws.Range("I6").Value = 1.047463333801878
$ws.Range("J6").Value = 1.055830233427535
$ws.Range("K6").Value = 1.065595296687436
$ws.Range("L6").Value = 1.061346224784179
$ws.Range("M6").Value = 1.07347464241321
$ws.Range("N6").Value = 1.022391333649569
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.051151126760102
$ws.Range("D7").Value = 1.062998896974177
$ws.Range("E7").Value = 1.058689674226953
$ws.Range("F7").Value = 1.070836677871157
$ws.Range("I7").Value = 1.047381185154253
$ws.Range("J7").Value = 1.055609320956365
$ws.Range("K7").Value = 1.065415776723261
$ws.Range("L7").Value = 1.061116906805137
$ws.Range("M7").Value = 1.07323495676384
$ws.Range("N7").Value = 1.022316292646489
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.049765522034182
$ws.Range("D8").Value = 1.062003774924669
$ws.Range("E8").Value = 1.057485690901674
$ws.Range("F8").Value = 1.069592429438406
$ws.Range("I8").Value = 1.047035748096836
$ws.Range("J8").Value = 1.054685540536524
$ws.Range("K8").Value = 1.064664552120447
$ws.Range("L8").Value = 1.060158488935134
$ws.Range("M8").Value = 1.072233263042069
$ws.Range("N8").Value = 1.022002235590374
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.047326978103081
$ws.Range("D9").Value = 1.060252166324506
$ws.Range("E9").Value = 1.055369371986689
$ws.Range("F9").Value = 1.067405361319272
$ws.Range("I9").Value = 1.046419270145882
$ws.Range("J9").Value = 1.053056240427605
$ws.Range("K9").Value = 1.063337608389703
$ws.Range("L9").Value = 1.058470032563073
$ws.Range("M9").Value = 1.07046878292002
$ws.Range("N9").Value = 1.021447345887628
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.045703358812151
$ws.Range("D10").Value = 1.059085787849819
$ws.Range("E10").Value = 1.053962075484328
$ws.Range("F10").Value = 1.06595104337968
$ws.Range("I10").Value = 1.046003081047127
$ws.Range("J10").Value = 1.051969075432343
$ws.Range("K10").Value = 1.062450881852459
$ws.Range("L10").Value = 1.057344701774134
$ws.Range("M10").Value = 1.069292936142362
$ws.Range("N10").Value = 1.021076440236779
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.045000808341842
$ws.Range("D11").Value = 1.058581069685001
$ws.Range("E11").Value = 1.053353560549833
$ws.Range("F11").Value = 1.065322203709399
$ws.Range("I11").Value = 1.045821636891893
$ws.Range("J11").Value = 1.051498098132819
$ws.Range("K11").Value = 1.062066431680166
$ws.Range("L11").Value = 1.056857501661216
$ws.Range("M11").Value = 1.06878390352643
$ws.Range("N11").Value = 1.020915605085531
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.04473992310823
$ws.Range("D12").Value = 1.058393645757991
$ws.Range("E12").Value = 1.053127660042138
$ws.Range("F12").Value = 1.065088759220292
$ws.Range("I12").Value = 1.045754055759768
$ws.Range("J12").Value = 1.05132312240649
$ws.Range("K12").Value = 1.061923556626085
$ws.Range("L12").Value = 1.05667654576665
$ws.Range("M12").Value = 1.068594844428948
$ws.Range("N12").Value = 1.020855829357526
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.044795880558686
$ws.Range("D13").Value = 1.058433846507145
$ws.Range("E13").Value = 1.053176110616498
$ws.Range("F13").Value = 1.065138827751712
$ws.Range("I13").Value = 1.045768560494503
$ws.Range("J13").Value = 1.051360656815922
$ws.Range("K13").Value = 1.061954207108233
$ws.Range("L13").Value = 1.056715360871724
$ws.Range("M13").Value = 1.068635397408104
$ws.Range("N13").Value = 1.020868653006424
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.044979241991504
$ws.Range("D14").Value = 1.058565576120704
$ws.Range("E14").Value = 1.053334884906709
$ws.Range("F14").Value = 1.065302904364771
$ws.Range("I14").Value = 1.045816054380219
$ws.Range("J14").Value = 1.051483635271503
$ws.Range("K14").Value = 1.062054623071457
$ws.Range("L14").Value = 1.056842543544052
$ws.Range("M14").Value = 1.068768275465287
$ws.Range("N14").Value = 1.020910664705357
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.045092226735969
$ws.Range("D15").Value = 1.058646745856607
$ws.Range("E15").Value = 1.053432728108283
$ws.Range("F15").Value = 1.06540401523957
$ws.Range("I15").Value = 1.045845292466842
$ws.Range("J15").Value = 1.051559401892288
$ws.Range("K15").Value = 1.062116482990452
$ws.Range("L15").Value = 1.056920906582638
$ws.Range("M15").Value = 1.068850148466187
$ws.Range("N15").Value = 1.020936544949796
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.045749995047472
$ws.Range("D16").Value = 1.059119291406925
$ws.Range("E16").Value = 1.05400247866928
$ws.Range("F16").Value = 1.06599279621811
$ws.Range("I16").Value = 1.046015096979283
$ws.Range("J16").Value = 1.052000327883333
$ws.Range("K16").Value = 1.062476386246713
$ws.Range("L16").Value = 1.057377037281336
$ws.Range("M16").Value = 1.069326721483965
$ws.Range("N16").Value = 1.021087109499571
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.046162726266205
$ws.Range("D17").Value = 1.059415796237875
$ws.Range("E17").Value = 1.0543600972556
$ws.Range("F17").Value = 1.066362361783706
$ws.Range("I17").Value = 1.046121281380155
$ws.Range("J17").Value = 1.052276848446794
$ws.Range("K17").Value = 1.062702012918895
$ws.Range("L17").Value = 1.057663176645235
$ws.Range("M17").Value = 1.069625694934398
$ws.Range("N17").Value = 1.02118149311141
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.046403512425199
$ws.Range("D18").Value = 1.059588774393253
$ws.Range("E18").Value = 1.05456877244826
$ws.Range("F18").Value = 1.066578008897312
$ws.Range("I18").Value = 1.046183098097776
$ws.Range("J18").Value = 1.052438116216753
$ws.Range("K18").Value = 1.062833569645335
$ws.Range("L18").Value = 1.057830084145725
$ws.Range("M18").Value = 1.069800092309882
$ws.Range("N18").Value = 1.02123652321339
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.046485622256758
$ws.Range("D19").Value = 1.059647760852644
$ws.Range("E19").Value = 1.054639939236129
$ws.Range("F19").Value = 1.066651553510881
$ws.Range("I19").Value = 1.04620415583252
$ws.Range("J19").Value = 1.052493100615088
$ws.Range("K19").Value = 1.062878419011161
$ws.Range("L19").Value = 1.057886996502629
$ws.Range("M19").Value = 1.069859559188929
$ws.Range("N19").Value = 1.021255283272537
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.046118439248212
$ws.Range("D20").Value = 1.059383980767791
$ws.Range("E20").Value = 1.054321719645234
$ws.Range("F20").Value = 1.066322702017641
$ws.Range("I20").Value = 1.046109901089355
$ws.Range("J20").Value = 1.052247182679398
$ws.Range("K20").Value = 1.062677810219483
$ws.Range("L20").Value = 1.057632475845351
$ws.Range("M20").Value = 1.069593616746737
$ws.Range("N20").Value = 1.021171368943955
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.044925244595917
$ws.Range("D21").Value = 1.058526783619788
$ws.Range("E21").Value = 1.053288126267459
$ws.Range("F21").Value = 1.065254584169782
$ws.Range("I21").Value = 1.045802073710433
$ws.Range("J21").Value = 1.051447422110793
$ws.Range("K21").Value = 1.062025055086275
$ws.Range("L21").Value = 1.056805091083356
$ws.Range("M21").Value = 1.068729145678809
$ws.Range("N21").Value = 1.020898294254528
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.044175459843427
$ws.Range("D22").Value = 1.057988124615279
$ws.Range("E22").Value = 1.052639012015254
$ws.Range("F22").Value = 1.064583795081069
$ws.Range("I22").Value = 1.045607462001577
$ws.Range("J22").Value = 1.050944385627664
$ws.Range("K22").Value = 1.061614219211086
$ws.Range("L22").Value = 1.056284950693815
$ws.Range("M22").Value = 1.068185724538745
$ws.Range("N22").Value = 1.020726402502006
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.044572894705162
$ws.Range("D23").Value = 1.058273649634708
$ws.Range("E23").Value = 1.052983048704773
$ws.Range("F23").Value = 1.064939318924023
$ws.Range("I23").Value = 1.045710730513717
$ws.Range("J23").Value = 1.051211073153672
$ws.Range("K23").Value = 1.061832050914243
$ws.Range("L23").Value = 1.056560680200709
$ws.Range("M23").Value = 1.068473792061689
$ws.Range("N23").Value = 1.020817544339134
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.046138450499927
$ws.Range("D24").Value = 1.059398356710752
$ws.Range("E24").Value = 1.05433906058108
$ws.Range("F24").Value = 1.066340622293819
$ws.Range("I24").Value = 1.046115043720131
$ws.Range("J24").Value = 1.052260587430678
$ws.Range("K24").Value = 1.062688746524345
$ws.Range("L24").Value = 1.057646348192918
$ws.Range("M24").Value = 1.069608111464311
$ws.Range("N24").Value = 1.021175943688231
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.047957035143124
$ws.Range("D25").Value = 1.060704764243319
$ws.Range("E25").Value = 1.055915863375341
$ws.Range("F25").Value = 1.067970118166731
$ws.Range("I25").Value = 1.04657956335664
$ws.Range("J25").Value = 1.053477625898942
$ws.Range("K25").Value = 1.063681028302164
$ws.Range("L25").Value = 1.058906487791622
$ws.Range("M25").Value = 1.070924863440945
$ws.Range("N25").Value = 1.021590971979206

Write-Host "Updated vm_pu results for 380 kV case"
